# Automatische test-sync: 2025-06-25 22:43:50
#
# Adds a second test-mail row to the "Logs" sheet, clears the stray
# "nan" placeholder that used to sit in E4, rolls the new category
# ("Openingstijden / Locatie") into the "Dashboard" summary sheet, and
# extends the conditional formatting + bar-chart ranges so the new row
# is picked up everywhere it needs to be.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs: row 4 had an un-filled "nan" placeholder reply -> clear it ---
$logs.Range("E4").Value = ""

# --- Logs: append the new test-mail as row 5 ---
$logs.Range("A5").Value = "Wanneer zijn jullie open?"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #2: Wanneer zijn jullie open?"
$logs.Range("D5").Value = "Openingstijden / Locatie"
$logs.Range("E5").Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F5").Value = "2025-06-25 22:43:26"
$logs.Range("G5").Value = "Ja"

# --- Logs: extend the two conditional-formatting blocks to cover row 5 ---
$catFormats = $logs.Range("D2:D4").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D5"))
}

$answeredFormats = $logs.Range("G2:G4").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G5"))
}

# --- Dashboard: add the new category's tally row ---
$dash.Range("A3").Value = "Openingstijden / Locatie"
$dash.Range("B3").Value = 1

# --- Dashboard chart: widen the category/value series references ---
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
